$wb = $excel.ActiveWorkbook
$excel.Left = -98
$excel.Top = -98
$excel.Width = 20715
$excel.Height = 13276
